# Apply the updated product-report layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: update the report date/time stamp ---
$ws.Range("D1").Value = "Fecha  2023-02-04 22:13:24"

# --- Row 3: remove old merges, add new header columns ---
$ws.Range("A3:B3").UnMerge()
$ws.Range("B3:C3").UnMerge()

$ws.Range("A3").Value = "Nombre"
$ws.Range("B3").Value = "Concentración"
$ws.Range("C3").Value = "Adicional"
$ws.Range("D3").Value = "Precio"
$ws.Range("E3").Value = "Tipo"
$ws.Range("F3").Value = "Presentacion"
$ws.Range("G3").Value = "Laboratorio"

# Give the newly added header cells the same (bold) formatting as the
# existing header cells, matching the row's look-and-feel.
$ws.Range("A3").Copy()
$ws.Range("C3:G3").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 4 data updates ---
$ws.Range("A4").Value = "buscapina "

# --- Row 5 data updates ---
$ws.Range("D5").Value = 5
$ws.Range("F5").Value = "Suero"

# --- Row 7 data updates ---
$ws.Range("F7").Value = "Suero"
